$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.286040782928467
$ws.Range("B1").Value = 1.97884476184845
$ws.Range("C1").Value = 5.507226467132568
$ws.Range("D1").Value = 1.920494794845581
$ws.Range("E1").Value = 1.10482931137085
